# Apply the "sertifikat, surat pernyataan, logo, fitur pencarian nama" data
# template edit: insert a new "Kota Lahir" (birth city) column between
# "Role" and "Tanggal Lahir", give the "password" column a real header,
# and replace/extend the sample rows accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

# ---------------------------------------------------------------------
# Header row (row 1) — shift "Tanggal Lahir" from D to E, add new "Kota
# Lahir" header in D, and give the previously-blank F1 a real "password"
# header (matching the bold style already used by the other headers).
# ---------------------------------------------------------------------
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial($xlPasteFormats)

$ws.Range("D1").Value = "Kota Lahir"
$ws.Range("E1").Value = "Tanggal Lahir"
$ws.Range("F1").Value = "password"

# ---------------------------------------------------------------------
# Row 2 — replace "Aqila Nur Azza" record with "Kemal S", add his birth
# city ("Surabaya") and a real date value for "Tanggal Lahir", and
# duplicate the username into the new "password" column.
# ---------------------------------------------------------------------
$ws.Range("A2").Value = "Kemal S"
$ws.Range("B2").Value = 2341760196

$ws.Range("C2").Copy()
$ws.Range("D2").PasteSpecial($xlPasteFormats)
$ws.Range("D2").Value = "Surabaya"

$ws.Range("D4").Copy()
$ws.Range("E2").PasteSpecial($xlPasteFormats)
$ws.Range("E2").Value = 37632

$ws.Range("C2").Copy()
$ws.Range("F2").PasteSpecial($xlPasteFormats)
$ws.Range("F2").Value = 2341760196

# ---------------------------------------------------------------------
# Row 3 — keep "Ismi Atika", add her birth city ("Bangkalan"), move her
# existing "15-01-2004" text into the "Tanggal Lahir" column, and
# duplicate her username into the new "password" column.
# ---------------------------------------------------------------------
$ws.Range("B3").Copy()
$ws.Range("F3").PasteSpecial($xlPasteFormats)
$ws.Range("F3").Value = 2341760036

$ws.Range("D3").Copy()
$ws.Range("E3").PasteSpecial($xlPasteFormats)
$ws.Range("E3").Value = "15-01-2004"

$ws.Range("C3").Copy()
$ws.Range("D3").PasteSpecial($xlPasteFormats)
$ws.Range("D3").Value = "Bangkalan"

$ws.Range("B3").Value = 2341760036

# ---------------------------------------------------------------------
# Column widths for the new layout.
# ---------------------------------------------------------------------
$ws.Columns("D").ColumnWidth = 12.666666666666666
$ws.Columns("E").ColumnWidth = 14.666666666666666
$ws.Columns("F").ColumnWidth = 16

# ---------------------------------------------------------------------
# Selection follows the edited "Tanggal Lahir" cell for row 3.
# ---------------------------------------------------------------------
$ws.Range("E3").Select()
